# Update cryptocurrency price/volume data to reflect refreshed GitHub Actions scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "260.97"
    "E2" = "1.74%"
    "D3" = "27.17"
    "E3" = "2.39%"
    "D4" = "4.689"
    "E4" = "0.48%"
    "E5" = "3.40%"
    "D6" = "6.657"
    "E6" = "0.74%"
    "D7" = "0.8534"
    "E7" = "-0.10%"
    "D8" = "0.9209"
    "E8" = "0.71%"
    "D9" = "0.1408"
    "E9" = "2.13%"
    "D10" = "0.04619"
    "E10" = "7.92%"
    "D11" = "0.07087"
    "E11" = "1.31%"
    "D12" = "0.03074"
    "E12" = "1.37%"
    "D13" = "0.09060"
    "E13" = "-0.55%"
    "D14" = "0.001528"
    "E14" = "0.35%"
    "D15" = "0.0006071"
    "E15" = "0.16%"
    "D16" = "0.006066"
    "E16" = "0.82%"
    "E17" = "-0.64%"
    "E18" = "0.33%"
    "E19" = "-0.20%"
    "E21" = "2.03%"
    "D22" = "4.102"
    "E22" = "5.72%"
    "E23" = "0.90%"
    "D24" = "0.001217"
    "E24" = "0.33%"
    "D25" = "0.003798"
    "E25" = "-18.43%"
    "E26" = "0.14%"
    "E27" = "3.51%"
    "D40" = "0.03873"
    "E40" = "2.05%"
    "D41" = "0.1113"
    "E41" = "1.25%"
    "D42" = "0.004080"
    "E42" = "-34.64%"
    "E43" = "12.39%"
    "E44" = "-3.99%"
    "D45" = "0.00005157"
    "E45" = "0.75%"
    "E46" = "0.14%"
    "E47" = "8.17%"
    "D48" = "0.1622"
    "E48" = "-32.64%"
    "E49" = "0.14%"
    "E50" = "0.14%"
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    # Ensure cell stays text-formatted (as it is stored as an inline string in the sheet)
    # so Excel doesn't reinterpret numeric-looking or percent-looking strings.
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
}

